# Change -999 ratings to 0 in Sheet1
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$cells = @("B2", "C2", "D3", "D6", "B8", "B9", "B10", "D11", "D12", "B15", "B16", "B17")
foreach ($cellRef in $cells) {
    $ws.Range($cellRef).Value = 0
}

# Update selection to match the saved cursor position in the diff
$ws.Range("E24").Select()

# Restore the window height captured at save time
$excel.ActiveWindow.Height = 17640
